# إضافة حدث جديد في Card24
# Fill the previously-blank trailing cells of row 19 with the literal "nan"
# placeholder (matching the convention used throughout the sheet for empty
# fields), then append a new service-event row (row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Row 19: fill blank cells B:K and P with "nan" ---
$ws.Range("B19").Value = "nan"
$ws.Range("C19").Value = "nan"
$ws.Range("D19").Value = "nan"
$ws.Range("E19").Value = "nan"
$ws.Range("F19").Value = "nan"
$ws.Range("G19").Value = "nan"
$ws.Range("H19").Value = "nan"
$ws.Range("I19").Value = "nan"
$ws.Range("J19").Value = "nan"
$ws.Range("K19").Value = "nan"
$ws.Range("P19").Value = "nan"

# --- Row 20: new service event ---
# Pre-format the card-number cell as Text so the numeric-looking value
# ("23") is stored as a string, matching the sheet's existing convention
# (every other cell in this table is text, e.g. A19 = "24").
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "23"
$ws.Range("L20").Value = "6\8\2024"
$ws.Range("N20").Value = "تم سن الفلاتس وتغيير اول جريده فوق المنشار"
$ws.Range("O20").Value = "الخبير"
